$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687816523147"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168785092949"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168785092949"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687851399388"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651168785236119"

# Sheet 1 (GNG_TO)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687816079183.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687816349216.csv"
$ws1.Range("B4").Value = "go_stims-1651168781635971.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168781651313.csv"

# Sheet 2 (NB_TO)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_7-16511687819235182.csv"
$ws2.Range("B3").Value = "OB-16511687832816553.csv"
$ws2.Range("B4").Value = "TB-1651168784443473.csv"
$ws2.Range("B5").Value = "ZB-match_5-16511687821657462.csv"
$ws2.Range("B6").Value = "TB-1651168783549834.csv"
$ws2.Range("B7").Value = "TB-1651168785079285.csv"
$ws2.Range("B8").Value = "ZB-match_5-16511687818978179.csv"
$ws2.Range("B9").Value = "OB-16511687827319396.csv"
$ws2.Range("B10").Value = "OB-16511687829064782.csv"

# Sheet 3 (RS_TO)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL_TO)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687851097288.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687850963194.csv"
$ws4.Range("B4").Value = "MM_stims-16511687851244512.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687851107218.csv"
$ws4.Range("B6").Value = "MM_stims-16511687851399388.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687851254554.csv"

# Sheet 5 (vSAT_TO)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687851448247.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511687852138677.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168785187522.csv"
$ws5.Range("B5").Value = "SAT_stims-1651168785171964.csv"
